$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Cells.Item(3, 6).Value = 3.5
$ws.Cells.Item(3, 7).Value = 4.9
$ws.Cells.Item(3, 8).Value = 1.88
$ws.Cells.Item(3, 9).Value = 2.3
$ws.Cells.Item(3, 10).Value = 3.05
$ws.Cells.Item(3, 11).Value = 3.95
$ws.Cells.Item(3, 12).Value = 1.41
$ws.Cells.Item(3, 13).Value = 1.07
$ws.Cells.Item(3, 14).Value = 3.2
$ws.Cells.Item(3, 15).Value = 1.33
$ws.Cells.Item(3, 16).Value = 1.84
$ws.Cells.Item(3, 17).Value = 1.94
$ws.Cells.Item(3, 18).Value = 1.32
$ws.Cells.Item(3, 19).Value = 3.15
$ws.Cells.Item(3, 20).Value = 1.79
$ws.Cells.Item(3, 21).Value = 2.02
$ws.Cells.Item(3, 22).Value = 1.83
$ws.Cells.Item(3, 24).Value = 16
$ws.Cells.Item(3, 25).Value = 10.5
$ws.Cells.Item(3, 26).Value = 15
$ws.Cells.Item(3, 27).Value = 29
$ws.Cells.Item(3, 28).Value = 17
$ws.Cells.Item(3, 29).Value = 9.4
$ws.Cells.Item(3, 30).Value = 12.5
$ws.Cells.Item(3, 31).Value = 26
$ws.Cells.Item(3, 32).Value = 36
$ws.Cells.Item(3, 33).Value = 19.5
$ws.Cells.Item(3, 34).Value = 22
$ws.Cells.Item(3, 35).Value = 44
$ws.Cells.Item(3, 36).Value = 100
$ws.Cells.Item(3, 37).Value = 60
$ws.Cells.Item(3, 38).Value = 70
$ws.Cells.Item(3, 39).Value = 120
$ws.Cells.Item(3, 40).Value = 70
$ws.Cells.Item(3, 41).Value = 19
$ws.Cells.Item(4, 6).Value = 3.2
$ws.Cells.Item(4, 7).Value = 4.2
$ws.Cells.Item(4, 8).Value = 2.14
$ws.Cells.Item(4, 9).Value = 2.66
$ws.Cells.Item(4, 11).Value = 3.7
$ws.Cells.Item(4, 12).Value = 1.38
$ws.Cells.Item(4, 13).Value = 1.08
$ws.Cells.Item(4, 14).Value = 3.2
$ws.Cells.Item(4, 15).Value = 1.37
$ws.Cells.Item(4, 16).Value = 1.74
$ws.Cells.Item(4, 17).Value = 2.08
$ws.Cells.Item(4, 18).Value = 1.28
$ws.Cells.Item(4, 19).Value = 3.8
$ws.Cells.Item(4, 20).Value = 1.81
$ws.Cells.Item(4, 22).Value = 1.66
$ws.Cells.Item(4, 23).Value = 1.31
$ws.Cells.Item(4, 24).Value = 15
$ws.Cells.Item(4, 25).Value = 11
$ws.Cells.Item(4, 26).Value = 17
$ws.Cells.Item(4, 27).Value = 970
$ws.Cells.Item(4, 28).Value = 15
$ws.Cells.Item(4, 29).Value = 9.199999999999999
$ws.Cells.Item(4, 30).Value = 13.5
$ws.Cells.Item(4, 31).Value = 970
$ws.Cells.Item(4, 32).Value = 970
$ws.Cells.Item(4, 33).Value = 18.5
$ws.Cells.Item(4, 34).Value = 23
$ws.Cells.Item(4, 35).Value = 55
$ws.Cells.Item(4, 36).Value = 90
$ws.Cells.Item(4, 37).Value = 60
$ws.Cells.Item(4, 38).Value = 75
$ws.Cells.Item(4, 39).Value = 140
$ws.Cells.Item(4, 40).Value = 65
$ws.Cells.Item(4, 41).Value = 970
$ws.Cells.Item(5, 6).Value = 1.94
$ws.Cells.Item(5, 7).Value = 2.14
$ws.Cells.Item(5, 8).Value = 3.9
$ws.Cells.Item(5, 10).Value = 3.35
$ws.Cells.Item(5, 11).Value = 3.9
$ws.Cells.Item(5, 12).Value = 1.41
$ws.Cells.Item(5, 13).Value = 1.07
$ws.Cells.Item(5, 14).Value = 3.4
$ws.Cells.Item(5, 15).Value = 1.33
$ws.Cells.Item(5, 16).Value = 1.83
$ws.Cells.Item(5, 17).Value = 1.99
$ws.Cells.Item(5, 18).Value = 1.32
$ws.Cells.Item(5, 19).Value = 3.55
$ws.Cells.Item(5, 20).Value = 1.82
$ws.Cells.Item(5, 21).Value = 2
$ws.Cells.Item(5, 23).Value = 1.87
$ws.Cells.Item(5, 24).Value = 16
$ws.Cells.Item(5, 25).Value = 17.5
$ws.Cells.Item(5, 26).Value = 34
$ws.Cells.Item(5, 27).Value = 110
$ws.Cells.Item(5, 28).Value = 10.5
$ws.Cells.Item(5, 29).Value = 8.6
$ws.Cells.Item(5, 30).Value = 21
$ws.Cells.Item(5, 31).Value = 60
$ws.Cells.Item(5, 32).Value = 15
$ws.Cells.Item(5, 33).Value = 13
$ws.Cells.Item(5, 34).Value = 23
$ws.Cells.Item(5, 35).Value = 80
$ws.Cells.Item(5, 36).Value = 30
$ws.Cells.Item(5, 37).Value = 24
$ws.Cells.Item(5, 38).Value = 48
$ws.Cells.Item(5, 39).Value = 130
$ws.Cells.Item(5, 40).Value = 19.5
$ws.Cells.Item(5, 41).Value = 65
$ws.Cells.Item(6, 6).Value = 1.12
$ws.Cells.Item(6, 7).Value = 1.3
$ws.Cells.Item(6, 8).Value = 10
$ws.Cells.Item(6, 10).Value = 4.3
$ws.Cells.Item(6, 11).Value = 10
$ws.Cells.Item(6, 23).Value = 4.3
$ws.Cells.Item(7, 6).Value = 1.08
$ws.Cells.Item(7, 9).Value = 16
$ws.Cells.Item(7, 10).Value = 1.07
$ws.Cells.Item(7, 13).Value = 1.07
$ws.Cells.Item(7, 15).Value = 1.07
$ws.Cells.Item(7, 17).Value = 1.07
$ws.Cells.Item(7, 19).Value = 1.07
$ws.Cells.Item(7, 22).Value = 1.06
$ws.Cells.Item(8, 17).Value = 1.66
$ws.Cells.Item(8, 22).Value = 1.05
$ws.Cells.Item(9, 6).Value = 1.57
$ws.Cells.Item(9, 7).Value = 1.64
$ws.Cells.Item(9, 11).Value = 4.1
$ws.Cells.Item(9, 15).Value = 1.47
$ws.Cells.Item(9, 17).Value = 2.4
$ws.Cells.Item(9, 19).Value = 4.7
$ws.Cells.Item(9, 21).Value = 1.66
$ws.Cells.Item(9, 23).Value = 2.56
$ws.Cells.Item(9, 26).Value = 75
$ws.Cells.Item(10, 6).Value = 3.45
$ws.Cells.Item(10, 7).Value = 3.75
$ws.Cells.Item(10, 8).Value = 2.08
$ws.Cells.Item(10, 9).Value = 2.18
$ws.Cells.Item(10, 11).Value = 4.2
$ws.Cells.Item(10, 12).Value = 1.24
$ws.Cells.Item(10, 21).Value = 2.6
$ws.Cells.Item(10, 22).Value = 1.85
$ws.Cells.Item(10, 23).Value = 1.37
$ws.Cells.Item(10, 32).Value = 32
$ws.Cells.Item(10, 33).Value = 18.5
$ws.Cells.Item(10, 36).Value = 60
$ws.Cells.Item(10, 37).Value = 38
$ws.Cells.Item(10, 38).Value = 40
$ws.Cells.Item(10, 39).Value = 60
$ws.Cells.Item(10, 40).Value = 26
$ws.Cells.Item(10, 41).Value = 11.5
$ws.Cells.Item(11, 14).Value = 4.4
$ws.Cells.Item(11, 15).Value = 1.28
$ws.Cells.Item(11, 16).Value = 2.16
$ws.Cells.Item(11, 17).Value = 1.82
$ws.Cells.Item(11, 18).Value = 1.45
$ws.Cells.Item(11, 19).Value = 3.1
$ws.Cells.Item(11, 20).Value = 1.78
$ws.Cells.Item(11, 21).Value = 2.22
$ws.Cells.Item(11, 24).Value = 16
$ws.Cells.Item(11, 27).Value = 20
$ws.Cells.Item(11, 38).Value = 60
$ws.Cells.Item(11, 39).Value = 90
$ws.Cells.Item(11, 40).Value = 50
$ws.Cells.Item(11, 41).Value = 11
$ws.Cells.Item(12, 6).Value = 3.55
$ws.Cells.Item(12, 8).Value = 2.16
$ws.Cells.Item(12, 9).Value = 2.2
$ws.Cells.Item(12, 23).Value = 1.38
$ws.Cells.Item(14, 9).Value = 3.55
$ws.Cells.Item(14, 11).Value = 3.95
$ws.Cells.Item(14, 13).Value = 1.06
$ws.Cells.Item(14, 14).Value = 3.75
$ws.Cells.Item(14, 15).Value = 1.29
$ws.Cells.Item(14, 17).Value = 1.87
$ws.Cells.Item(14, 18).Value = 1.37
$ws.Cells.Item(14, 19).Value = 3.15
$ws.Cells.Item(14, 20).Value = 1.71
$ws.Cells.Item(14, 21).Value = 2.14
$ws.Cells.Item(14, 24).Value = 970
$ws.Cells.Item(14, 28).Value = 970
$ws.Cells.Item(14, 29).Value = 8.4
$ws.Cells.Item(14, 32).Value = 970
$ws.Cells.Item(14, 35).Value = 970
$ws.Cells.Item(14, 37).Value = 970
$ws.Cells.Item(15, 7).Value = 2.7
$ws.Cells.Item(15, 10).Value = 3.65
$ws.Cells.Item(16, 8).Value = 3.9
$ws.Cells.Item(16, 20).Value = 2.32
$ws.Cells.Item(16, 28).Value = 7.2
$ws.Cells.Item(16, 29).Value = 8.6
$ws.Cells.Item(16, 33).Value = 15.5
$ws.Cells.Item(17, 7).Value = 3.85
$ws.Cells.Item(17, 8).Value = 2.46
$ws.Cells.Item(17, 9).Value = 2.72
$ws.Cells.Item(17, 10).Value = 2.86
$ws.Cells.Item(17, 11).Value = 2.98
$ws.Cells.Item(17, 12).Value = 1.57
$ws.Cells.Item(17, 22).Value = 1.58
$ws.Cells.Item(17, 23).Value = 1.35
$ws.Cells.Item(17, 28).Value = 970
$ws.Cells.Item(19, 7).Value = 2.24
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1.73
$ws.Cells.Item(19, 15).Value = 1.76
$ws.Cells.Item(19, 16).Value = 1.35
$ws.Cells.Item(19, 17).Value = 3.35
$ws.Cells.Item(19, 19).Value = 7.6
$ws.Cells.Item(19, 23).Value = 1.8
$ws.Cells.Item(19, 24).Value = 6.4
$ws.Cells.Item(20, 22).Value = 1.41
$ws.Cells.Item(22, 7).Value = 2.18
$ws.Cells.Item(22, 8).Value = 4.6
$ws.Cells.Item(22, 16).Value = 1.57
$ws.Cells.Item(22, 21).Value = 1.78
$ws.Cells.Item(22, 23).Value = 1.84
$ws.Cells.Item(23, 6).Value = 2.38
$ws.Cells.Item(23, 7).Value = 2.58
$ws.Cells.Item(23, 10).Value = 2.86
$ws.Cells.Item(23, 13).Value = 1.14
$ws.Cells.Item(23, 22).Value = 1.31
$ws.Cells.Item(24, 7).Value = 3.95
$ws.Cells.Item(24, 9).Value = 2.64
$ws.Cells.Item(24, 10).Value = 2.86
$ws.Cells.Item(24, 12).Value = 1.62
$ws.Cells.Item(24, 14).Value = 2.44
$ws.Cells.Item(24, 15).Value = 1.59
$ws.Cells.Item(24, 16).Value = 1.48
$ws.Cells.Item(24, 17).Value = 2.8
$ws.Cells.Item(24, 18).Value = 1.17
$ws.Cells.Item(24, 19).Value = 6
$ws.Cells.Item(24, 20).Value = 2.18
$ws.Cells.Item(24, 21).Value = 1.72
$ws.Cells.Item(24, 22).Value = 1.61
$ws.Cells.Item(24, 23).Value = 1.34
$ws.Cells.Item(24, 24).Value = 9.4
$ws.Cells.Item(24, 39).Value = 260
$ws.Cells.Item(24, 40).Value = 120
$ws.Cells.Item(25, 17).Value = 2.46
